# Weekly crime-data refresh for the 47th Precinct CompStat report
# (Volume 32, Number 49 -> 50; report week 12/1-12/7/2025 -> 12/8-12/14/2025),
# including the new week's complaint counts / percent changes and the
# addition of the "***.*" masked-value placeholder string used for rows
# where the percent-change cannot be computed (previously suppressed via a
# literal "0" placeholder instead).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header: volume/issue number and the week-ending date range ---
$ws.Range("A8").Value = "Volume 32   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/8/2025  Through  12/14/2025"

# --- Cells whose content flips from a real number to the masked-value
#     placeholder text ("0" or "***.*"). These must stay text cells (not
#     numeric 0), so we nudge Excel into text mode with a leading
#     apostrophe, then restore the original numeric-column formatting by
#     pasting formats only from a cell that already carries that style. ---
$maskedCells = @(
    @("D14", "0"),
    @("E14", "***.*"),
    @("D22", "0"),
    @("E22", "***.*"),
    @("D31", "0"),
    @("E31", "***.*"),
    @("C33", "0"),
    @("D33", "0"),
    @("E33", "***.*")
)
foreach ($cell in $maskedCells) {
    $target = $cell[0]
    $text = $cell[1]
    $ws.Range($target).Value = "'" + $text
    $ws.Range("C14").Copy()
    $ws.Range($target).PasteSpecial(-4122)
}

# --- Updated weekly / 28-day / YTD / 2-year complaint counts and percent
#     changes for the crime-category table (rows 14-33) and the historical
#     perspective table (rows 40-46 only shift which label they point at;
#     no values change there). ---
$updatedValues = @(
    @("N14", -64.705882352941),
    @("C15", 2),
    @("D15", 8),
    @("E15", -75),
    @("F15", 5),
    @("G15", 11),
    @("H15", -54.545454545454),
    @("I15", 60),
    @("J15", 50),
    @("K15", 20),
    @("L15", 46.341463414634),
    @("M15", 57.894736842105),
    @("N15", -21.052631578947),
    @("C16", 6),
    @("D16", 11),
    @("E16", -45.454545454545),
    @("F16", 27),
    @("G16", 45),
    @("H16", -40),
    @("I16", 464),
    @("J16", 523),
    @("K16", -11.281070745697),
    @("L16", 2.202643171806),
    @("M16", 10.213776722090),
    @("N16", -65.450483991064),
    @("C17", 22),
    @("D17", 17),
    @("E17", 29.411764705882),
    @("F17", 82),
    @("G17", 57),
    @("H17", 43.859649122807),
    @("I17", 977),
    @("J17", 777),
    @("K17", 25.740025740025),
    @("L17", 16.033254156769),
    @("M17", 134.292565947242),
    @("N17", 18.567961165048),
    @("C18", 5),
    @("D18", 4),
    @("E18", 25),
    @("F18", 15),
    @("G18", 19),
    @("H18", -21.052631578947),
    @("I18", 241),
    @("J18", 252),
    @("K18", -4.365079365079),
    @("L18", -3.6),
    @("M18", -33.972602739726),
    @("N18", -86.528787031861),
    @("C19", 16),
    @("D19", 16),
    @("E19", 0),
    @("F19", 76),
    @("G19", 67),
    @("H19", 13.432835820895),
    @("I19", 963),
    @("J19", 853),
    @("K19", 12.895662368112),
    @("L19", 25.390625),
    @("M19", 225.337837837838),
    @("N19", 82.732447817836),
    @("C20", 9),
    @("D20", 7),
    @("E20", 28.571428571428),
    @("F20", 42),
    @("G20", 38),
    @("H20", 10.526315789473),
    @("I20", 604),
    @("J20", 523),
    @("K20", 15.487571701720),
    @("L20", 1.003344481605),
    @("M20", 72.571428571428),
    @("N20", -61.057382333978),
    @("C21", 60),
    @("D21", 63),
    @("E21", -4.761904761904),
    @("F21", 247),
    @("G21", 239),
    @("H21", 3.347280334728),
    @("I21", 3321),
    @("J21", 2986),
    @("K21", 11.219022103148),
    @("L21", 12.082348970637),
    @("M21", 73.965426925091),
    @("N21", -45.947265625),
    @("F22", 2),
    @("G22", 3),
    @("H22", -33.333333333333),
    @("I22", 24),
    @("J22", 30),
    @("K22", -20),
    @("L22", -11.111111111111),
    @("M22", -17.241379310344),
    @("C23", 2),
    @("D23", 8),
    @("E23", -75),
    @("F23", 13),
    @("G23", 15),
    @("H23", -13.333333333333),
    @("I23", 158),
    @("J23", 123),
    @("K23", 28.455284552845),
    @("L23", 61.224489795918),
    @("M23", 119.444444444444),
    @("C24", 22),
    @("D24", 34),
    @("E24", -35.294117647058),
    @("F24", 124),
    @("G24", 126),
    @("H24", -1.587301587301),
    @("I24", 1415),
    @("J24", 1287),
    @("K24", 9.945609945609),
    @("L24", 11.154752553024),
    @("M24", 92.255434782608),
    @("C25", 1),
    @("D25", 9),
    @("E25", -88.888888888888),
    @("F25", 19),
    @("G25", 29),
    @("H25", -34.482758620689),
    @("I25", 209),
    @("J25", 378),
    @("K25", -44.708994708994),
    @("L25", -40.114613180515),
    @("C26", 24),
    @("D26", 17),
    @("E26", 41.176470588235),
    @("F26", 103),
    @("G26", 72),
    @("H26", 43.055555555555),
    @("I26", 1174),
    @("J26", 1060),
    @("K26", 10.754716981132),
    @("L26", 22.546972860125),
    @("M26", 26.100966702470),
    @("C27", 2),
    @("D27", 8),
    @("E27", -75),
    @("F27", 6),
    @("G27", 11),
    @("H27", -45.454545454545),
    @("I27", 77),
    @("J27", 70),
    @("K27", 10),
    @("L27", 16.666666666666),
    @("C28", 2),
    @("D28", 1),
    @("E28", 100),
    @("F28", 6),
    @("G28", 5),
    @("H28", 20),
    @("I28", 103),
    @("J28", 80),
    @("K28", 28.75),
    @("L28", 37.333333333333),
    @("C29", 1),
    @("D29", 1),
    @("E29", 0),
    @("F29", 4),
    @("G29", 5),
    @("H29", -20),
    @("I29", 42),
    @("J29", 36),
    @("K29", 16.666666666666),
    @("L29", 10.526315789473),
    @("M29", -34.375),
    @("N29", -66.929133858267),
    @("C30", 1),
    @("D30", 1),
    @("E30", 0),
    @("F30", 4),
    @("G30", 4),
    @("H30", 0),
    @("I30", 31),
    @("J30", 26),
    @("K30", 19.230769230769),
    @("L30", -6.060606060606),
    @("M30", -40.384615384615),
    @("N30", -74.166666666666),
    @("G33", 2),
    @("H33", -50)
)
foreach ($cell in $updatedValues) {
    $ws.Range($cell[0]).Value = $cell[1]
}
